{"js": "// 3rd commit of pinkesh\n// The document ends with an (empty) paragraph that only holds the\n// \"_GoBack\" bookmark left over from the previous edit session. This\n// commit types \"Commit pinkesh 3\" into that paragraph (ahead of the\n// bookmark, exactly like Word does when you resume typing at the\n// remembered cursor position) and then presses Enter, which leaves a\n// brand new empty paragraph behind it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The paragraph we want is the very last paragraph in the document.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Type the new commit line at the start of that (currently empty)\n// paragraph, keeping the trailing bookmark where it was.\nlastParagraph.insertText(\"Commit pinkesh 3\", Word.InsertLocation.start);\n\n// Hit Enter: a fresh, empty paragraph follows.\nlastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# 3rd commit of pinkesh\n# The document ends with an (empty) paragraph that only holds the\n# \"_GoBack\" bookmark left over from the previous edit session. This\n# commit types \"Commit pinkesh 3\" into that paragraph (ahead of the\n# bookmark, exactly like Word does when you resume typing at the\n# remembered cursor position) and then presses Enter, which leaves a\n# brand new empty paragraph behind it.\n\n$d = $word.ActiveDocument\n\n# The paragraph we want is the very last paragraph in the document.\n$lastParagraph = $d.Paragraphs.Last\n\n# Type the new commit line at the start of that (currently empty)\n# paragraph, keeping the trailing bookmark where it was.\n$lastParagraph.Range.InsertBefore(\"Commit pinkesh 3\")\n\n# Hit Enter: a fresh, empty paragraph follows.\n$lastParagraph.Range.InsertParagraphAfter()\n"}
